# update pages, also closes #8
#
# Adds a new column M ("fin_vyuct_verejne") to Sheet1, filling in its
# header (styled like the other header cells) and the 25 data rows below
# it. Also refreshes the dimension/used range implicitly by writing into
# M1:M26, and nudges the one cached literal that shifted by a float ULP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell M1, matching style of the existing header row ---
$ws.Range("M1").Value = "fin_vyuct_verejne"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108   # xlCenter

# --- New data column M2:M26 ---
$mValues = @{
    2  = 11553428901.86
    3  = 7259831475.81
    4  = 3318112356.43
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1706372594.94
    10 = 18561765.15
    11 = 29750632143.26718
    12 = 637662428.21
    13 = 8458455415.152819
    14 = 0
    15 = 2205418812.46
    16 = 113866716.77
    17 = 775373909.08
    18 = 0
    19 = 23761275973.88
    20 = 3479432893.66
    21 = 2201763329.6
    22 = 2578418745.04
    23 = 67715138.94
    24 = 0
    25 = 0
    26 = 0
}

foreach ($row in $mValues.Keys) {
    $ws.Cells.Item($row, 13).Value = $mValues[$row]
}

# --- Tiny recalculated-literal drift on L13 (same magnitude, last ULP) ---
$ws.Range("L13").Value = [double]"-6.374381599011212e-08"

Write-Output "Column M populated; dimension now $($ws.UsedRange.Address())"
